$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 12.99892624393442
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 26.973800054782714
$ws.Range("E2").Value = 27.480116759324574

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 19.04410753028111
$ws.Range("D3").Value = 14.132875048404799
$ws.Range("E3").Value = 31.265119917477023

# Update selection to match target
$ws.Range("B1:E3").Select()
